# Auto-generated Excel COM-interop script
# Applies numeric corrections to specific cells across multiple sheets
# as described by the commit diff (scheduled runner profit recalculations).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 4835.718
$ws.Range("J17").Value = 4884.0264
$ws.Range("L17").Value = 14652.0792
$ws.Range("N17").Value = -14988.0792

# Row 132
$ws.Range("H132").Value = 1471
$ws.Range("I132").Value = 1444.68
$ws.Range("J132").Value = 1800
$ws.Range("K132").Value = 4334.04
$ws.Range("L132").Value = 5400
$ws.Range("M132").Value = -1804.04
$ws.Range("N132").Value = -10460

$ws = $wb.Worksheets.Item("ARM")
# Row 22
$ws.Range("H22").Value = 4998
$ws.Range("I22").Value = 4998
$ws.Range("K22").Value = 4998
$ws.Range("M22").Value = -4699

# Row 28
$ws.Range("H28").Value = 87499.5
$ws.Range("I28").Value = 60000
$ws.Range("K28").Value = 60000
$ws.Range("M28").Value = -59808

# Row 32
$ws.Range("H32").Value = 25373240
$ws.Range("I32").Value = 28295964
$ws.Range("K32").Value = 28295964
$ws.Range("M32").Value = -28295677

# Row 45
$ws.Range("H45").Value = 3546.5454
$ws.Range("I45").Value = 2649.75
$ws.Range("J45").Value = 4059
$ws.Range("K45").Value = 2649.75
$ws.Range("L45").Value = 4059
$ws.Range("M45").Value = -2272.75
$ws.Range("N45").Value = -4813

# Row 61
$ws.Range("H61").Value = 2939.24
$ws.Range("I61").Value = 2604.5
$ws.Range("J61").Value = 3800
$ws.Range("K61").Value = 2604.5
$ws.Range("L61").Value = 3800
$ws.Range("M61").Value = -2392.5
$ws.Range("N61").Value = -4224

# Row 74
$ws.Range("H74").Value = 2689.1177
$ws.Range("I74").Value = 2431.923
$ws.Range("J74").Value = 3525
$ws.Range("K74").Value = 2431.923
$ws.Range("L74").Value = 3525
$ws.Range("M74").Value = -1557.923
$ws.Range("N74").Value = -5273

# Row 77
$ws.Range("H77").Value = 2689.1177
$ws.Range("I77").Value = 2431.923
$ws.Range("J77").Value = 3525
$ws.Range("K77").Value = 12159.615
$ws.Range("L77").Value = 17625
$ws.Range("M77").Value = -7791.614999999998
$ws.Range("N77").Value = -26361

# Row 99
$ws.Range("H99").Value = 87499.5
$ws.Range("I99").Value = 60000
$ws.Range("K99").Value = 60000
$ws.Range("M99").Value = -57005

# Row 102
$ws.Range("H102").Value = 1623.68
$ws.Range("I102").Value = 1412.7727
$ws.Range("K102").Value = 1412.7727
$ws.Range("M102").Value = 209.2273

# Row 108
$ws.Range("H108").Value = 93195.60000000001
$ws.Range("J108").Value = 93195.60000000001
$ws.Range("L108").Value = 93195.60000000001
$ws.Range("N108").Value = -100875.6

# Row 124
$ws.Range("H124").Value = 80000
$ws.Range("J124").Value = 80000
$ws.Range("L124").Value = 80000
$ws.Range("N124").Value = -89820

# Row 132
$ws.Range("H132").Value = 289061.94
$ws.Range("I132").Value = 387217.7
$ws.Range("J132").Value = 5500.8887
$ws.Range("K132").Value = 1161653.1
$ws.Range("L132").Value = 16502.6661
$ws.Range("M132").Value = -1159123.1
$ws.Range("N132").Value = -21562.6661

# Row 136
$ws.Range("H136").Value = 2939.24
$ws.Range("I136").Value = 2604.5
$ws.Range("J136").Value = 3800
$ws.Range("K136").Value = 7813.5
$ws.Range("L136").Value = 11400
$ws.Range("M136").Value = -5263.5
$ws.Range("N136").Value = -16500

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 2718.5
$ws.Range("I105").Value = 2278.7856
$ws.Range("J105").Value = 3334.1
$ws.Range("K105").Value = 2278.7856
$ws.Range("L105").Value = 3334.1
$ws.Range("M105").Value = -531.7856000000002
$ws.Range("N105").Value = -6828.1

# Row 126
$ws.Range("H126").Value = 169000
$ws.Range("J126").Value = 169000
$ws.Range("L126").Value = 169000
$ws.Range("N126").Value = -178880

# Row 134
$ws.Range("H134").Value = 4467753.5
$ws.Range("I134").Value = 5104968.5
$ws.Range("J134").Value = 7250
$ws.Range("K134").Value = 15314905.5
$ws.Range("L134").Value = 21750
$ws.Range("M134").Value = -15312370.5
$ws.Range("N134").Value = -26820

$ws = $wb.Worksheets.Item("CRP")
# Row 64
$ws.Range("H64").Value = 46995
$ws.Range("J64").Value = 46995
$ws.Range("L64").Value = 46995
$ws.Range("N64").Value = -47491

# Row 67
$ws.Range("H67").Value = 46995
$ws.Range("J67").Value = 46995
$ws.Range("L67").Value = 46995
$ws.Range("N67").Value = -48711

# Row 86
$ws.Range("H86").Value = 34984.52
$ws.Range("I86").Value = 22311.062
$ws.Range("J86").Value = 40617.168
$ws.Range("K86").Value = 22311.062
$ws.Range("L86").Value = 40617.168
$ws.Range("M86").Value = -21188.062
$ws.Range("N86").Value = -42863.168

# Row 89
$ws.Range("H89").Value = 34984.52
$ws.Range("I89").Value = 22311.062
$ws.Range("J89").Value = 40617.168
$ws.Range("K89").Value = 111555.31
$ws.Range("L89").Value = 203085.84
$ws.Range("M89").Value = -105939.31
$ws.Range("N89").Value = -214317.84

# Row 100
$ws.Range("H100").Value = 60005.89
$ws.Range("J100").Value = 60005.89
$ws.Range("L100").Value = 60005.89
$ws.Range("N100").Value = -62169.89

# Row 132
$ws.Range("H132").Value = 5015.6665
$ws.Range("I132").Value = 4960.1113
$ws.Range("J132").Value = 5182.3335
$ws.Range("K132").Value = 14880.3339
$ws.Range("L132").Value = 15547.0005
$ws.Range("M132").Value = -12350.3339
$ws.Range("N132").Value = -20607.0005

# Row 134
$ws.Range("H134").Value = 2618.4546
$ws.Range("I134").Value = 2422.5557
$ws.Range("K134").Value = 7267.6671
$ws.Range("M134").Value = -4732.6671

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2152.2104
$ws.Range("I102").Value = 2010.5714
$ws.Range("K102").Value = 2010.5714
$ws.Range("M102").Value = -388.5714

# Row 132
$ws.Range("H132").Value = 2312.8572
$ws.Range("I132").Value = 1781.6666
$ws.Range("J132").Value = 5500
$ws.Range("K132").Value = 5344.9998
$ws.Range("L132").Value = 16500
$ws.Range("M132").Value = -2814.9998
$ws.Range("N132").Value = -21560

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1245.3572
$ws.Range("I22").Value = 1351.625
$ws.Range("J22").Value = 1103.6666
$ws.Range("K22").Value = 1351.625
$ws.Range("L22").Value = 1103.6666
$ws.Range("M22").Value = -1056.625
$ws.Range("N22").Value = -1693.6666

# Row 24
$ws.Range("H24").Value = 14862.5
$ws.Range("I24").Value = 14862.5
$ws.Range("K24").Value = 14862.5
$ws.Range("M24").Value = -14519.5

# Row 27
$ws.Range("H27").Value = 1245.3572
$ws.Range("I27").Value = 1351.625
$ws.Range("J27").Value = 1103.6666
$ws.Range("K27").Value = 1351.625
$ws.Range("L27").Value = 1103.6666
$ws.Range("M27").Value = -1244.625
$ws.Range("N27").Value = -1317.6666

# Row 55
$ws.Range("H55").Value = 475.91666
$ws.Range("I55").Value = 422
$ws.Range("K55").Value = 422
$ws.Range("M55").Value = -249

# Row 133
$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -55060

# Row 136
$ws.Range("H136").Value = 11567.5
$ws.Range("I136").Value = 5916.364
$ws.Range("K136").Value = 17749.092
$ws.Range("M136").Value = -15199.092

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 36236.8
$ws.Range("I132").Value = 41248.348
$ws.Range("K132").Value = 123745.044
$ws.Range("M132").Value = -121215.044
